# Rotate the species/record data among rows 25, 26 and 27.
# Row 25 gets what used to be in row 26 (pre-edit),
# Row 26 gets what used to be in row 27 (pre-edit),
# Row 27 gets what used to be in row 25 (pre-edit).
# Columns C, D, I, P, S, T, U, V, W, Y, Z, AA, AB, AD, AE, AG, AT, AW, AX, AY
# are identical across the three rows and are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A25").Value = 112111665
$ws.Range("B25").Value = 78578
$ws.Range("E25").Value = 6458
$ws.Range("F25").Value = "Lunglav"
$ws.Range("G25").Value = "Lobaria pulmonaria"
$ws.Range("H25").Value = "(L.) Hoffm."
$ws.Range("Q25").Value = 553817.5168682858
$ws.Range("R25").Value = 6956613.349017856

$ws.Range("A26").Value = 112111666
$ws.Range("B26").Value = 78578
$ws.Range("E26").Value = 6458
$ws.Range("F26").Value = "Lunglav"
$ws.Range("G26").Value = "Lobaria pulmonaria"
$ws.Range("H26").Value = "(L.) Hoffm."
$ws.Range("Q26").Value = 553943.5691689024
$ws.Range("R26").Value = 6956459.040529874

$ws.Range("A27").Value = 112111680
$ws.Range("B27").Value = 89423
$ws.Range("E27").Value = 5432
$ws.Range("F27").Value = "Granticka"
$ws.Range("G27").Value = "Porodaedalea chrysoloma"
$ws.Range("H27").Value = "(Fr.) Fiasson & Niemelä"
$ws.Range("Q27").Value = 554111.3423843421
$ws.Range("R27").Value = 6956617.245754472
